$d = $word.ActiveDocument

# Remove the leading empty paragraphs that only hold a horizontal-rule
# drawing (w:pict/v:rect) in front of the "Introduction" bookmark/heading.
# Walk from the top of the document and delete any paragraph whose text
# is empty (i.e. just the paragraph mark) until we reach real content.
while ($d.Paragraphs.Count -gt 0 -and $d.Paragraphs(1).Range.Text.Trim().Length -eq 0) {
    $d.Paragraphs(1).Range.Delete()
}

Write-Output $d.Paragraphs.Count
